$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the very top; this shifts all existing rows (1-9) down to (2-10)
$ws.Rows.Item(1).Insert()

# New column J header (row 1): "loop unrolling" label with a yellow highlight fill
$ws.Range("J1").Value2 = "loop unrolling"
$ws.Range("J1").Interior.Color = 65535

# New column J sub-header (row 2): same "Kernel time (ms)" label used in column H
$ws.Range("J2").Value2 = "Kernel time (ms)"

# New column J data values (rows 3-10); row 7 stays empty, matching the gap in column E-H
$ws.Range("J3").Value2 = 2817.8130000000001
$ws.Range("J4").Value2 = 461.50299999999999
$ws.Range("J5").Value2 = 374.96
$ws.Range("J6").Value2 = 380.839
$ws.Range("J8").Value2 = 398.89299999999997
$ws.Range("J9").Value2 = 377.12
$ws.Range("J10").Value2 = 431.24099999999999

# Give column J the same width as column H
$ws.Columns.Item(10).ColumnWidth = 13.14

# Update the active selection to match the saved workbook state
$ws.Range("E17").Select()
